# Swap the "detalhar diária" step (TC2, row 18) with the
# "cancelar diária" step (TC4, row 32), so that the step order becomes:
#   TC2 -> cancelar diária
#   TC3 -> analisar prestação de contas (unchanged)
#   TC4 -> detalhar diária
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$stepDetalhar = "Beneficiário Clica em detalhar diária."
$resultDetalhar = "SYSTEM Apresenta a tela de Detalhar Diárias"

$stepCancelar = "Beneficiário Clica em cancelar diária."
$resultCancelar = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"

# Row 18 (TC2 step) currently holds the "detalhar" text; put "cancelar" there.
$ws.Range("B18").Value = $stepCancelar
$ws.Range("D18").Value = $resultCancelar

# Row 32 (TC4 step) currently holds the "cancelar" text; put "detalhar" there.
$ws.Range("B32").Value = $stepDetalhar
$ws.Range("D32").Value = $resultDetalhar
